# Updates Results_Production_xgb_Solina_wm.xlsx ("Production_Predictions" sheet)
# with the latest production-prediction batch:
#   - Date (col A) / Interval (col B) / Prediction (col C) values shift to the
#     newer forecast run for rows 2-68.
#   - Rows 69-97 only carry a refreshed Prediction value; the Date/Interval
#     columns for those trailing rows are cleared (no longer populated for
#     this export).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry is "row,col,value" (col 1=A/Date, 2=B/Interval, 3=C/Prediction).
$edits = "2,1,45337;2,2,10;2,3,0.3129999935626984;3,1,45337;3,2,11;3,3,0.3930000066757202;4,1,45337;4,2,12;4,3,0.425000011920929;5,1,45337;5,2,13;5,3,0.4020000100135803;6,1,45337;6,2,14;6,3,0.3459999859333038;7,1,45337;7,2,15;7,3,0.2029999941587448;8,1,45337;8,2,16;8,3,0.1620000004768372;9,1,45337;9,2,17;9,3,0.01400000043213367;10,1,45337;10,2,18;11,2,19;12,2,20;13,2,21;14,2,22;15,2,23;16,1,45338;16,2,0;17,1,45338;17,2,1;18,1,45338;18,2,2;18,3,0.001000000047497451;19,1,45338;19,2,3;19,3,0.001000000047497451;20,1,45338;20,2,4;20,3,0.001000000047497451;21,1,45338;21,2,5;21,3,0.001000000047497451;22,1,45338;22,2,6;22,3,0.001000000047497451;23,1,45338;23,2,7;23,3,0.01099999994039536;24,1,45338;24,2,8;24,3,0.08600000292062759;25,1,45338;25,2,9;25,3,0.1850000023841858;26,1,45338;26,2,10;26,3,0.3310000002384186;27,1,45338;27,2,11;27,3,0.4099999964237213;28,1,45338;28,2,12;28,3,0.4079999923706055;29,1,45338;29,2,13;29,3,0.3919999897480011;30,1,45338;30,2,14;30,3,0.300000011920929;31,1,45338;31,2,15;31,3,0.1979999989271164;32,1,45338;32,2,16;32,3,0.1389999985694885;33,1,45338;33,2,17;33,3,0.02099999971687794;34,1,45338;34,2,18;35,2,19;36,2,20;37,2,21;38,2,22;39,2,23;40,1,45339;40,2,0;41,1,45339;41,2,1;42,1,45339;42,2,2;42,3,0.001000000047497451;43,1,45339;43,2,3;43,3,0.001000000047497451;44,1,45339;44,2,4;44,3,0.001000000047497451;45,1,45339;45,2,5;45,3,0.001000000047497451;46,1,45339;46,2,6;46,3,0.001000000047497451;47,1,45339;47,2,7;47,3,0.01099999994039536;48,1,45339;48,2,8;48,3,0.08100000023841858;49,1,45339;49,2,9;49,3,0.1940000057220459;50,3,0.1979999989271164;51,3,0.1389999985694885;52,3,0.02099999971687794;66,3,0.01099999994039536;68,3,0.1940000057220459;69,3,0.2569999992847443;70,3,0.3720000088214874;71,3,0.414000004529953;72,3,0.4050000011920929;73,3,0.3120000064373016;74,3,0.2689999938011169;75,3,0.239999994635582;76,3,0.1389999985694885;77,3,0.05900000035762787;78,3,0.04199999943375587;79,3,0.0390000008046627;80,3,0.03400000184774399;81,3,0.03400000184774399;82,3,0.03400000184774399;83,3,0.05799999833106995;84,3,0.05900000035762787;85,3,0.05799999833106995;86,3,0.05700000002980232;87,3,0.05700000002980232;88,3,0.05700000002980232;89,3,0.05700000002980232;90,3,0.07400000095367432;91,3,0.1570000052452087;92,3,0.1580000072717667;93,3,0.1759999990463257;94,3,0.210999995470047;95,3,0.239999994635582;96,3,0.257999986410141;97,3,0.2310000061988831"

foreach ($entry in $edits.Split(";")) {
    $parts = $entry.Split(",")
    $r = [int]$parts[0]
    $c = [int]$parts[1]
    $v = [double]$parts[2]
    $ws.Cells.Item($r, $c).Value = $v
}

# The trailing rows (69-97) no longer report a Date/Interval for this export
# run - only the refreshed Prediction (col C) remains.
$ws.Range("A69:B97").Clear()
